$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Example" mini-table in columns J:Q (rows 1-8) + notes (rows 10-13)
# ---------------------------------------------------------------------------

# --- Row 1: banner "Example" across J1:Q1 -----------------------------------
$ws.Range("J1").Value = "Example"
$ws.Range("J1:Q1").Style = "Accent6"
$ws.Range("J1:Q1").Merge()
$ws.Range("J1:Q1").HorizontalAlignment = -4108   # xlCenter

# --- Row 2: column headers ---------------------------------------------------
$ws.Range("J2:Q2").Style = "20% - Accent6"
$ws.Range("J2").Value = "Group No."
$ws.Range("K2").Value = "Project name"
$ws.Range("L2").Value = "Company"
$ws.Range("M2").Value = "Student ID"
$ws.Range("N2").Value = "Name"
$ws.Range("O2").Value = "Role"
$ws.Range("P2").Value = "Advisor"
$ws.Range("Q2").Value = "Co-advisor"

# --- Rows 3-8: two example groups, base style first --------------------------
$ws.Range("J3:Q8").Style = "20% - Accent6"

# Group 1 (rows 3-5)
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = "Rocket to the moon"
$ws.Range("L3").Value = "Alphabet Inc."
$ws.Range("M3").Value = 65130500201
$ws.Range("M4").Value = 65130500202
$ws.Range("N3").Value = "Nadech Kugimiya"
$ws.Range("N4").Value = "Peem Wasu"
$ws.Range("O3").Value = "Tester"
$ws.Range("O4").Value = "UX/UI"
$ws.Range("P3").Value = "Chonlameth Arpnikanondt"
$ws.Range("Q3").Value = "Tuul Triyason"

$ws.Range("J3:J5").Merge()
$ws.Range("K3:K5").Merge()
$ws.Range("L3:L5").Merge()
$ws.Range("P3:P5").Merge()
$ws.Range("Q3:Q5").Merge()

# Group 2 (rows 6-8)
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = "Rocket to Mars"
$ws.Range("L6").Value = "Microsoft Corporation"
$ws.Range("M6").Value = 65130500203
$ws.Range("N6").Value = "Janis Star"
$ws.Range("O6").Value = "Frontend developer"
$ws.Range("P6").Value = "Vithida Chongsuphajaisiddhi"

$ws.Range("J6:J8").Merge()
$ws.Range("K6:K8").Merge()
$ws.Range("L6:L8").Merge()
$ws.Range("P6:P8").Merge()
$ws.Range("Q6:Q8").Merge()

# Alignment: Group No. column -> left/vcenter ; text columns -> centered
$ws.Range("J3:J8").HorizontalAlignment = -4131   # xlLeft
$ws.Range("J3:J8").VerticalAlignment = -4108     # xlCenter

$ws.Range("K3:L8").HorizontalAlignment = -4108
$ws.Range("P3:Q8").HorizontalAlignment = -4108

# --- Rows 10-13: notes --------------------------------------------------------
$ws.Range("J10:P13").Style = "20% - Accent2"
$ws.Range("J10").Value = "You can find the lecturer name at https://www.sit.kmutt.ac.th/lecturer/"
$ws.Range("J11").Value = "Group member limit to 3 person."
$ws.Range("J12").Value = "Complete all columns except Co-advisor, which is optional."
$ws.Range("J13").Value = "Please re-check the correctness of the information."

# ---------------------------------------------------------------------------
# Column widths for the new J:Q columns
# ---------------------------------------------------------------------------
$ws.Columns("J").ColumnWidth = 9.166666666666666
$ws.Columns("K").ColumnWidth = 17.166666666666668
$ws.Columns("L").ColumnWidth = 19.736979166666668
$ws.Columns("M").ColumnWidth = 11.166666666666666
$ws.Columns("N").ColumnWidth = 15.736979166666666
$ws.Columns("O").ColumnWidth = 17.307291666666668
$ws.Columns("P").ColumnWidth = 26.451822916666668
$ws.Columns("Q").ColumnWidth = 12.022135416666666

# ---------------------------------------------------------------------------
# Selection (matches the final saved cursor position)
# ---------------------------------------------------------------------------
[void]$ws.Range("H5:H7").Select()

Write-Output "edit complete"
